$wb = $excel.ActiveWorkbook

# Add a new worksheet for Orders, placed after "Client List"
$clientList = $wb.Worksheets.Item("Client List")
$ordersSheet = $wb.Worksheets.Add([System.Type]::Missing, $clientList)
$ordersSheet.Name = "Orders"

# Header row
$ordersSheet.Range("A1").Value = "Item"
$ordersSheet.Range("B1").Value = "Quantity"
$ordersSheet.Range("C1").Value = "Customer ID"
$ordersSheet.Range("D1").Value = "Order ID"

# Data rows
$ordersSheet.Range("A2").Value = "yacht club soda"
$ordersSheet.Range("B2").Value = 4
$ordersSheet.Range("C2").Value = 706167
$ordersSheet.Range("D2").Value = 66146698

$ordersSheet.Range("A3").Value = "quohogs"
$ordersSheet.Range("B3").Value = 4
$ordersSheet.Range("C3").Value = 706167
$ordersSheet.Range("D3").Value = 93349827

# Autofit columns A (Item) and C (Customer ID) to fit their widest entries
$ordersSheet.Columns.Item(1).AutoFit() | Out-Null
$ordersSheet.Columns.Item(3).AutoFit() | Out-Null
$ordersSheet.Columns.Item(1).ColumnWidth = 28.35
$ordersSheet.Columns.Item(3).ColumnWidth = 10.35

# The Customer ID column on Client List also gets an auto-fit width now
# that Customer ID values are referenced/used by the new Orders sheet
$clientList.Columns.Item(9).AutoFit() | Out-Null
$clientList.Columns.Item(9).ColumnWidth = 10.35

# Make Orders the active sheet/tab
$ordersSheet.Activate()
$ordersSheet.Range("A2").Select() | Out-Null
